$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $TextValue)
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $TextValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.997.87"
Set-TextValue "D3" "1.910.58"
Set-TextValue "E3" "  +0.87%  "
Set-TextValue "D4" "1.000"
Set-TextValue "D5" "0.8050"
Set-TextValue "E5" "  +5.60%  "
Set-TextValue "D6" "242.12"
Set-TextValue "D7" "0.9996"
Set-TextValue "E7" "  +0.04%  "
Set-TextValue "D8" "0.3160"
Set-TextValue "E8" "  +3.86%  "
Set-TextValue "D9" "26.39"
Set-TextValue "E9" "  +3.84%  "
Set-TextValue "D10" "0.06912"
Set-TextValue "E10" "  +1.40%  "
Set-TextValue "D11" "0.07995"
Set-TextValue "E11" "  +0.04%  "
Set-TextValue "D12" "1.920.61"
Set-TextValue "E12" "  +1.48%  "
Set-TextValue "D13" "0.7400"
Set-TextValue "E13" "  -1.09%  "
Set-TextValue "D14" "5.207"
Set-TextValue "E14" "  +0.03%  "
Set-TextValue "D15" "93.08"
Set-TextValue "E15" "  +2.28%  "
Set-TextValue "D16" "30.002.00"
Set-TextValue "E16" "  +0.60%  "
Set-TextValue "E17" "  +0.92%  "
Set-TextValue "D18" "5.885"
Set-TextValue "D19" "245.99"
Set-TextValue "E19" "  +4.75%  "
Set-TextValue "D20" "0.000007738"
Set-TextValue "E20" "  +0.86%  "
Set-TextValue "D21" "0.9994"
Set-TextValue "E21" "  +0.00%  "
Set-TextValue "D22" "2.153.82"
Set-TextValue "E22" "  +1.13%  "
Set-TextValue "D23" "1.000"
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "6.855"
Set-TextValue "E24" "  -1.25%  "
Set-TextValue "E25" "  +1.43%  "
Set-TextValue "D26" "9.238"
Set-TextValue "E26" "  +0.12%  "
Set-TextValue "D27" "0.1426"
Set-TextValue "D29" "2.037"
Set-TextValue "E29" "  -0.34%  "
Set-TextValue "D30" "1.362"
Set-TextValue "E30" "  +1.47%  "
Set-TextValue "D31" "1.515"
Set-TextValue "E31" "  +0.32%  "
Set-TextValue "D32" "4.314"
Set-TextValue "E32" "  +0.91%  "
Set-TextValue "E33" "  +1.66%  "
Set-TextValue "D34" "0.05486"
Set-TextValue "E34" "  +2.51%  "
Set-TextValue "D35" "1.264"
Set-TextValue "E35" "  +1.44%  "
Set-TextValue "D36" "0.7349"
Set-TextValue "E36" "  +1.14%  "
Set-TextValue "D37" "2.719"
Set-TextValue "E37" "  +0.27%  "
Set-TextValue "D38" "0.01928"
Set-TextValue "E38" "  +0.26%  "
Set-TextValue "D39" "2.791"
Set-TextValue "E39" "  +0.66%  "
Set-TextValue "D40" "6.181"
Set-TextValue "E40" "  +0.17%  "
Set-TextValue "D41" "0.4426"
Set-TextValue "E41" "  +0.55%  "
Set-TextValue "D42" "72.45"
Set-TextValue "E42" "  +0.49%  "
Set-TextValue "D43" "0.9993"
Set-TextValue "E43" "  -0.06%  "
Set-TextValue "D44" "0.8373"
Set-TextValue "E44" "  +1.77%  "
Set-TextValue "D45" "1.879"
Set-TextValue "E45" "  -1.70%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D46" "100.62"
Set-TextValue "E46" "  -0.37%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.558"
Set-TextValue "E47" "  -0.28%  "
Set-TextValue "D48" "9.763"
Set-TextValue "E48" "  +0.09%  "
Set-TextValue "D49" "984.72"
Set-TextValue "E49" "  +6.53%  "
Set-TextValue "D50" "2.060.90"
Set-TextValue "E50" "  +1.03%  "
Set-TextValue "D51" "36.27"
Set-TextValue "E51" "  +0.28%  "
